$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46061 -> 46062) for every data row (rows 2 through 278).
$ws.Range("C2:C278").Value = 46062
